# Update the FATEcol crosstab in the single table of the document.
# Row 16 = "  Voyage completed as intended"
#   French column (col 5): 34 -> 35
#   Total column  (col 7): 302 -> 303
# Row 17 = "  Original goal thwarted before disembarking slaves"
#   French column (col 5): 3 -> 2
#   Total column  (col 7): 7 -> 6
# (One French-flagged voyage was reclassified from "thwarted before
# disembarking slaves" to "completed as intended".)

$d = $word.ActiveDocument
$t = $d.Tables.Item(1)

$t.Cell(16, 5).Range.Text = "35"
$t.Cell(16, 7).Range.Text = "303"

$t.Cell(17, 5).Range.Text = "2"
$t.Cell(17, 7).Range.Text = "6"
